# Update "想去人数" (interest count) figures on the "展览" and "全部类型" sheets
# to reflect the latest scrape (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 530
$ws1.Range("F7").Value = 1602
$ws1.Range("F10").Value = 1351
$ws1.Range("F12").Value = 17
$ws1.Range("F13").Value = 235
$ws1.Range("F17").Value = 8
$ws1.Range("F18").Value = 246
$ws1.Range("F21").Value = 189

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 530
$ws4.Range("F7").Value = 1602
$ws4.Range("F11").Value = 1351
$ws4.Range("F13").Value = 17
$ws4.Range("F14").Value = 235
$ws4.Range("F18").Value = 8
$ws4.Range("F19").Value = 246
$ws4.Range("F22").Value = 189
